# "Added Validation Details" - append two newly-validated customer rows
# (Test User / Test User2) to the onboarding sheet, mark the contact
# number column as text (so the leading zero is preserved), turn the
# e-mail address of the first new row into a mailto: hyperlink, size the
# columns to fit their content and leave the selection on I10, matching
# what Excel does when a user types this data in by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Test User ---------------------------------------------------
$ws.Range("A12").Value = "Test User"
$ws.Range("B12").Value = "20 Symonds Street, Auckland, NZ"
$ws.Range("D12").Value = "test.user@email.com"
$ws.Range("E12").Value = "Test"
$ws.Range("F12").Value = "User"

# Turn the e-mail address into a clickable mailto hyperlink (adds the
# "Hyperlink" cell style / font automatically).
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:test.user@email.com")

# Store the phone number as text so the leading zero survives.
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "02041166935"

# --- Row 13: Test User2 ---------------------------------------------------
$ws.Range("A13").Value = "Test User2"
$ws.Range("B13").Value = "21 Symonds Street, Auckland, NZ"
$ws.Range("E13").Value = "Test"
$ws.Range("F13").Value = "User2"

$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "02041166935"

# --- Column sizing to fit the new, longer content -------------------------
$ws.Columns.Item(1).ColumnWidth = 19.15
$ws.Columns.Item(2).ColumnWidth = 49.65
$ws.Columns.Item(3).ColumnWidth = 10.15
$ws.Columns.Item(4).ColumnWidth = 24.65

# --- Leave the selection where the author left it --------------------------
$ws.Range("I10").Select()
